$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("26÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("55÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷8=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("61÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("95÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷2=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("78÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷3=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("15÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("51÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("41÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("46÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("20÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷8=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("57÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("14÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷8=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("11÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷4=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("14÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷6=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("23÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("79÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷3=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷9=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("18÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("20÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷3=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("39÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷8=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷4=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("77÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("30÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷4=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("54÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷7=", 2) | Out-Null
$r = $d.Content
$r.Find.Execute("93÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "43÷2=", 2) | Out-Null
$r = $d.Content
